# 14-Apr-2024: Administrator functions were implemented.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current row 8 ("test paper" / "testpaper.pdf"),
# pushing the existing rows 8-9 down to rows 10-11.
$ws.Rows("8:9").Insert()

# New rows 8-9: candidates data/candidates.xlsx, test results/scores.xlsx.
# Filled in the filenames first, then the labels, to mirror the author's
# original editing order (and resulting shared-string table order).
$ws.Range("B8").Value = "candidates.xlsx"
$ws.Range("B9").Value = "scores.xlsx"
$ws.Range("A9").Value = "test results"
$ws.Range("A8").Value = "candidates data"

# Match the formatting used by the row above (A7:B7) for the two new rows.
$ws.Range("A8:B9").Font.Size = $ws.Range("A7").Font.Size
$ws.Range("B8:B9").HorizontalAlignment = $ws.Range("B7").HorizontalAlignment
$ws.Rows("8:9").RowHeight = $ws.Rows("7").RowHeight

# Update the selection to match the recorded UI state.
$ws.Range("A8").Select()

$wb.Save()
